$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together with each record when rows are re-keyed.
$cols = @("A","B","D","E","F","G","H","Q","R")

# Snapshot current ("before") values for every row that will be touched,
# since the update is a cyclic re-shuffle of rows and later writes must
# not clobber values still needed as a source for another row.
$rows = @(2,3,4,6,7,8,9,10,13)
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# Destination row -> source row (i.e. new row content = old content of source row).
$mapping = @{
    2  = 13
    3  = 6
    4  = 10
    6  = 8
    7  = 4
    8  = 9
    9  = 7
    10 = 2
    13 = 3
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value = $srcVals[$c]
    }
}
